# "Generate Report for Handoff"
#
# The localization-status report is regenerated: the zh-cn language row has
# moved from "In Translation" to "Ready for handoff" (the de-de row was
# already reporting that status) and the "Latest Handoff / HO Xliff
# Generate" timestamps are bumped forward a few seconds to the new
# generation run. The wider status text also means the report's datetime
# columns are redrawn a bit wider on every sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Excel's ColumnWidth setter quantizes to whole screen pixels (width in
# characters -> round(width*MDW + 5px padding) -> back to characters, with
# MDW = 6px for the workbook's default Calibri 11 font). To land the stored
# column width as close as possible to a desired "characters" value, solve
# that rounding for the input that reproduces the target pixel count.
function Set-PreciseColumnWidth($range, [double]$targetWidth) {
    $mdw = 6.0
    $desiredPixels = [Math]::Round($targetWidth * $mdw)
    $range.ColumnWidth = ($desiredPixels - 5) / $mdw
}

# --- Overview sheet: zh-cn/de-de status cells + the HO Xliff generate date ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-29 06:59:32"

Set-PreciseColumnWidth $overview.Columns.Item(5) 17.2159881591797
Set-PreciseColumnWidth $overview.Columns.Item(6) 17.2159881591797

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-29 06:59:28"

Set-PreciseColumnWidth $zhcn.Columns.Item(3) 17.2159881591797

# --- de-de sheet: Status + Latest Handoff Datetime ---
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-29 06:59:32"

Set-PreciseColumnWidth $dede.Columns.Item(3) 17.2159881591797
